# Add "Create New User" (testSuccessfulAdminAddAttachment fix) and
# "CreateNewUser_Test" (testSuccessfulNewUserCreate) test data blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCaseDataSets")

# -----------------------------------------------------------------
# 1. Fix pre-existing "AddAttachment" block (rows 27-28): the header
#    label in B27 was mistakenly the literal "Richmond County" value
#    instead of the column name "client", and the password/case
#    number test values in row 28 were stale.
# -----------------------------------------------------------------
$ws.Range("B27").Value = "client"
$ws.Range("D28").Value = "password1"
$ws.Range("E28").Value = "071-4-088-00-0"

# -----------------------------------------------------------------
# 2. New "CreateNewUser_Test" block: header row 31, data row 32, and
#    a blank bordered row 33 underneath (matching the look of the
#    other test blocks on this sheet).
# -----------------------------------------------------------------

# Clone the look of the existing header/data/blank rows (27/28/--)
# onto the new rows 31/32/33 before writing any values.
$ws.Range("A27:H27").Copy()
$ws.Range("A31:L31").PasteSpecial(-4122)

$ws.Range("A28:H28").Copy()
$ws.Range("A32:L32").PasteSpecial(-4122)
$ws.Range("A28:H28").Copy()
$ws.Range("A33:L33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row
$ws.Range("A31").Value = "testSuccessfulNewUserCreate"
$ws.Range("B31").Value = "username"
$ws.Range("C31").Value = "password"
$ws.Range("D31").Value = "nusername"
$ws.Range("E31").Value = "status"
$ws.Range("F31").Value = "role"
$ws.Range("G31").Value = "telenumber"
$ws.Range("H31").Value = "email"
$ws.Range("I31").Value = "npassword"
$ws.Range("J31").Value = "confirmpw"
$ws.Range("K31").Value = "path"
$ws.Range("L31").Value = "imgname"

# Data row
$ws.Range("A32").Value = "testSuccessfulNewUserCreate"
$ws.Range("B32").Value = "Jason Lee"
$ws.Range("C32").Value = "password1"
$ws.Range("D32").Value = "NewClientUser"
$ws.Range("E32").Value = "Active"
$ws.Range("F32").Value = "Client"
$ws.Range("G32").Value = "770.330.1111"
$ws.Range("H32").Value = "newclient@user2.com"
$ws.Range("I32").Value = "clientuser1"
$ws.Range("J32").Value = "clientuser1"
$ws.Range("K32").Value = "C:\\testfolder\\"
$ws.Range("L32").Value = "toshiya3.jpg"

# H32 keeps the "no explicit style" look it has in the original file.
$ws.Range("H32").Style = "Normal"
$ws.Range("H32").Value = "newclient@user2.com"

# -----------------------------------------------------------------
# 3. Extend the trailing blank filler rows (34-35 already existed,
#    36-41 are brand new) across columns A and I:L so the whole block
#    lines up with the new 12-column width.
# -----------------------------------------------------------------
$ws.Range("A34:L41").NumberFormat = "@"

# -----------------------------------------------------------------
# 4. Column widths for the three new columns (I already existed).
# -----------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 19.166666666666668
$ws.Columns.Item(10).ColumnWidth = 14
$ws.Columns.Item(11).ColumnWidth = 14
$ws.Columns.Item(12).ColumnWidth = 11.5

# -----------------------------------------------------------------
# 5. View state: scrolled down a bit with L32 selected, matching
#    where the new data was entered.
# -----------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("L32").Select()
